$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @("group11","group10","group13","group12","group14","group5","group4","group7","group6","group1","group3","group2","group9","group8")
$values = @(31,81,53,56,96,56,98,80,86,37,58,32,51,89)

for ($i = 0; $i -lt $groups.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $groups[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("D8").Select()
